$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 443.75
$ws.Range("I5").Value = 350
$ws.Range("J5").Value = 725
$ws.Range("K5").Value = 350
$ws.Range("L5").Value = 725
$ws.Range("M5").Value = -235
$ws.Range("N5").Value = -955
# Row 17
$ws.Range("H17").Value = 850.3692
$ws.Range("J17").Value = 850.3692
$ws.Range("L17").Value = 2551.1076
$ws.Range("N17").Value = -2887.1076
# Row 19
$ws.Range("H19").Value = 3442.889
$ws.Range("I19").Value = 1498.5
$ws.Range("J19").Value = 3998.4285
$ws.Range("K19").Value = 1498.5
$ws.Range("L19").Value = 3998.4285
$ws.Range("M19").Value = -1323.5
$ws.Range("N19").Value = -4348.4285
# Row 129
$ws.Range("H129").Value = 900.98114
$ws.Range("J129").Value = 877.4902
$ws.Range("L129").Value = 2632.4706
$ws.Range("N129").Value = -12632.4706
# Row 137
$ws.Range("H137").Value = 1292.0968
$ws.Range("I137").Value = 1217.5
$ws.Range("J137").Value = 1680
$ws.Range("K137").Value = 3652.5
$ws.Range("L137").Value = 5040
$ws.Range("M137").Value = -1102.5
$ws.Range("N137").Value = -10140
# Row 138
$ws.Range("H138").Value = 2742.8462
$ws.Range("J138").Value = 3957.5557
$ws.Range("L138").Value = 11872.6671
$ws.Range("N138").Value = -22152.6671

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2599.305
$ws.Range("I32").Value = 1845.1389
$ws.Range("J32").Value = 8029.3
$ws.Range("K32").Value = 1845.1389
$ws.Range("L32").Value = 8029.3
$ws.Range("M32").Value = -1558.1389
$ws.Range("N32").Value = -8603.299999999999
# Row 61
$ws.Range("I61").Value = 3406
$ws.Range("J61").Value = 7749.75
$ws.Range("K61").Value = 3406
$ws.Range("L61").Value = 7749.75
$ws.Range("M61").Value = -3194
$ws.Range("N61").Value = -8173.75
# Row 74
$ws.Range("H74").Value = 1866.079
$ws.Range("I74").Value = 1790.909
$ws.Range("K74").Value = 1790.909
$ws.Range("M74").Value = -916.9090000000001
# Row 77
$ws.Range("H77").Value = 1866.079
$ws.Range("I77").Value = 1790.909
$ws.Range("K77").Value = 8954.545
$ws.Range("M77").Value = -4586.545
# Row 110
$ws.Range("H110").Value = 1603.2667
$ws.Range("I110").Value = 876.9167
$ws.Range("J110").Value = 4508.6665
$ws.Range("K110").Value = 876.9167
$ws.Range("L110").Value = 4508.6665
$ws.Range("M110").Value = 1168.0833
$ws.Range("N110").Value = -8598.666499999999
# Row 122
$ws.Range("H122").Value = 1275
$ws.Range("I122").Value = 1275
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3825
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1375
$ws.Range("N122").ClearContents() | Out-Null
# Row 132
$ws.Range("H132").Value = 3023.6155
$ws.Range("I132").Value = 3177.25
$ws.Range("J132").Value = 2955.3333
$ws.Range("K132").Value = 9531.75
$ws.Range("L132").Value = 8865.999899999999
$ws.Range("M132").Value = -7001.75
$ws.Range("N132").Value = -13925.9999
# Row 136
$ws.Range("I136").Value = 3406
$ws.Range("J136").Value = 7749.75
$ws.Range("K136").Value = 10218
$ws.Range("L136").Value = 23249.25
$ws.Range("M136").Value = -7668
$ws.Range("N136").Value = -28349.25
# Row 139
$ws.Range("H139").Value = 49000
$ws.Range("J139").Value = 49000
$ws.Range("L139").Value = 49000
$ws.Range("N139").Value = -59280

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 57
$ws.Range("H57").Value = 49700
$ws.Range("J57").Value = 49700
$ws.Range("L57").Value = 49700
$ws.Range("N57").Value = -51140
# Row 94
$ws.Range("H94").Value = 298
$ws.Range("I94").Value = 308.2069
$ws.Range("K94").Value = 308.2069
$ws.Range("M94").Value = 142.7931
# Row 134
$ws.Range("H134").Value = 9310.056
$ws.Range("I134").Value = 9310.056
$ws.Range("K134").Value = 27930.168
$ws.Range("M134").Value = -25395.168
# Row 136
$ws.Range("H136").Value = 49700
$ws.Range("J136").Value = 49700
$ws.Range("L136").Value = 49700
$ws.Range("N136").Value = -59900
# Row 137
$ws.Range("H137").Value = 61375
$ws.Range("J137").Value = 61375
$ws.Range("L137").Value = 61375
$ws.Range("N137").Value = -71575

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1919.5714
$ws.Range("I31").Value = 1197.9412
$ws.Range("K31").Value = 1197.9412
$ws.Range("M31").Value = -902.9412
# Row 34
$ws.Range("H34").Value = 1919.5714
$ws.Range("I34").Value = 1197.9412
$ws.Range("K34").Value = 1197.9412
$ws.Range("M34").Value = -995.9412
# Row 107
$ws.Range("H107").Value = 327.10715
$ws.Range("I107").Value = 337.2
$ws.Range("J107").Value = 301.875
$ws.Range("K107").Value = 337.2
$ws.Range("L107").Value = 301.875
$ws.Range("M107").Value = 1582.8
$ws.Range("N107").Value = -4141.875
# Row 122
$ws.Range("H122").Value = 1637.8334
$ws.Range("J122").Value = 700
$ws.Range("L122").Value = 2100
$ws.Range("N122").Value = -7000
# Row 132
$ws.Range("H132").Value = 2248.348
$ws.Range("I132").Value = 1177
$ws.Range("J132").Value = 3641.1
$ws.Range("K132").Value = 3531
$ws.Range("L132").Value = 10923.3
$ws.Range("M132").Value = -1001
$ws.Range("N132").Value = -15983.3
# Row 134
$ws.Range("H134").Value = 965
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents() | Out-Null

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 43
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents() | Out-Null
# Row 114
$ws.Range("H114").Value = 3078.25
$ws.Range("J114").Value = 3883
$ws.Range("L114").Value = 11649
$ws.Range("N114").Value = -18157
# Row 129
$ws.Range("H129").Value = 28759.424
$ws.Range("I129").Value = 460.625
$ws.Range("J129").Value = 41336.668
$ws.Range("K129").Value = 1381.875
$ws.Range("L129").Value = 124010.004
$ws.Range("M129").Value = 3618.125
$ws.Range("N129").Value = -134010.004
# Row 131
$ws.Range("H131").Value = 764.8200000000001
$ws.Range("J131").Value = 780.7234
$ws.Range("L131").Value = 2342.1702
$ws.Range("N131").Value = -12422.1702
# Row 132
$ws.Range("H132").Value = 1200
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents() | Out-Null

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 3156.3635
$ws.Range("I102").Value = 3703.5
$ws.Range("K102").Value = 3703.5
$ws.Range("M102").Value = -2081.5
# Row 126
$ws.Range("H126").Value = 1826847.6
$ws.Range("I126").Value = 2224779
$ws.Range("J126").Value = 168800
$ws.Range("K126").Value = 6674337
$ws.Range("L126").Value = 506400
$ws.Range("M126").Value = -6671867
$ws.Range("N126").Value = -511340
# Row 132
$ws.Range("H132").Value = 6413389
$ws.Range("I132").Value = 12822112
$ws.Range("K132").Value = 38466336
$ws.Range("M132").Value = -38463806

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2839.5264
$ws.Range("I16").Value = 2827.0833
$ws.Range("K16").Value = 2827.0833
$ws.Range("M16").Value = -2657.0833
# Row 40
$ws.Range("H40").Value = 7359.8
$ws.Range("I40").Value = 2266.3333
$ws.Range("K40").Value = 2266.3333
$ws.Range("M40").Value = -2130.3333
# Row 132
$ws.Range("H132").Value = 1738.4474
$ws.Range("I132").Value = 1502.0625
$ws.Range("K132").Value = 4506.1875
$ws.Range("M132").Value = -1976.1875
# Row 136
$ws.Range("H136").Value = 3185.95
$ws.Range("J136").Value = 5437.375
$ws.Range("L136").Value = 16312.125
$ws.Range("N136").Value = -21412.125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 651.3333
$ws.Range("I100").Value = 419.2
$ws.Range("J100").Value = 941.5
$ws.Range("K100").Value = 838.4
$ws.Range("L100").Value = 1883
$ws.Range("M100").Value = -297.4
$ws.Range("N100").Value = -2965
# Row 132
$ws.Range("H132").Value = 1026.3214
$ws.Range("I132").Value = 778.38464
$ws.Range("K132").Value = 2335.15392
$ws.Range("M132").Value = 194.8460800000003
# Row 136
$ws.Range("H136").Value = 2939.4783
$ws.Range("I136").Value = 4261.5
$ws.Range("J136").Value = 2234.4
$ws.Range("K136").Value = 12784.5
$ws.Range("L136").Value = 6703.200000000001
$ws.Range("M136").Value = -10234.5
$ws.Range("N136").Value = -11803.2
